$wb = $excel.ActiveWorkbook

# 1. Rename sheet "baseline-forecasts" to "external-forecasts"
$wsForecasts = $wb.Worksheets.Item("baseline-forecasts")
$wsForecasts.Name = "external-forecasts"

# 2. Add trailing inflation variables (dns1, dns2, dns3) to the
#    baseline-variables table, rows 23-25, column A (varname)
$wsVars = $wb.Worksheets.Item("baseline-variables")
$wsVars.Range("A23").Value = "dns1"
$wsVars.Range("A24").Value = "dns2"
$wsVars.Range("A25").Value = "dns3"

# 3. Update the active selection on baseline-variables to B23
$wsVars.Activate()
$wsVars.Range("B23").Select()
